$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cell that already carries the workbook's plain "Century 12"
# font (font id 2 in the original styles.xml) with no extra color/theme
# info attached. Copying format from it (instead of poking
# Font.Name/Font.Size directly) lets COM reuse/derive the existing font
# table entries instead of minting a brand new "Calibri+theme" font for
# every cell we touch.
$fontTemplate = $ws.Cells.Item(7, 1)

function SetCell($row, $col, $val, $style) {
    $cell = $ws.Range($ws.Cells.Item($row, $col), $ws.Cells.Item($row, $col))
    if ($val -ne $null) {
        $cell.Value = $val
    }
    if ($style -ne "") {
        $fontTemplate.Copy()
        $cell.PasteSpecial(-4122)
        $cell.HorizontalAlignment = -4108

        if ($style -eq "bold") {
            $cell.Font.Bold = $true
        } elseif ($style -eq "red") {
            $cell.Font.Color = 255
        } elseif ($style -eq "green") {
            $cell.Font.Color = 32768
        } elseif ($style -eq "blue") {
            $cell.Font.Color = 16711680
        }

        $cell.Borders.LineStyle = 1
    }
}

SetCell 8 1 $null ""
SetCell 8 2 $null ""
SetCell 8 3 $null ""
SetCell 8 4 $null ""
SetCell 8 5 $null ""
SetCell 9 1 $null "plain"
SetCell 9 2 'Right' "plain"
SetCell 9 3 'Wrong' "plain"
SetCell 9 4 'Not Attempt' "plain"
SetCell 9 5 'Max' "plain"
SetCell 10 1 'No.' "plain"
SetCell 10 2 7 "green"
SetCell 10 3 11 "red"
SetCell 10 4 10 "plain"
SetCell 10 5 28 "plain"
SetCell 11 1 'Marking' "plain"
SetCell 11 2 5 "green"
SetCell 11 3 -1 "red"
SetCell 11 4 0 "plain"
SetCell 11 5 $null "plain"
SetCell 12 1 'Total' "plain"
SetCell 12 2 35 "green"
SetCell 12 3 -11 "red"
SetCell 12 4 $null "plain"
SetCell 12 5 '24/140' "blue"
SetCell 13 1 $null ""
SetCell 13 2 $null ""
SetCell 13 3 $null ""
SetCell 13 4 $null ""
SetCell 13 5 $null ""
SetCell 14 1 $null ""
SetCell 14 2 $null ""
SetCell 14 3 $null ""
SetCell 14 4 $null ""
SetCell 14 5 $null ""
SetCell 15 1 'Student Ans' "bold"
SetCell 15 2 'Correct Ans' "bold"
SetCell 15 4 'Student Ans' "bold"
SetCell 15 5 'Correct Ans' "bold"
SetCell 16 1 'Option D' "red"
SetCell 16 2 'Option A' "blue"
SetCell 16 4 $null "red"
SetCell 16 5 'Option A' "blue"
SetCell 17 1 'Option D' "green"
SetCell 17 2 'Option D' "blue"
SetCell 17 4 $null "red"
SetCell 17 5 'Option C' "blue"
SetCell 18 1 'Option A' "red"
SetCell 18 2 'Option B' "blue"
SetCell 18 4 'Option B' "red"
SetCell 18 5 'Option D' "blue"
SetCell 19 1 'Option C' "green"
SetCell 19 2 'Option C' "blue"
SetCell 20 1 'Option A' "red"
SetCell 20 2 'Option B' "blue"
SetCell 21 1 'Option C' "green"
SetCell 21 2 'Option C' "blue"
SetCell 22 1 $null "red"
SetCell 22 2 'Option D' "blue"
SetCell 23 1 $null "red"
SetCell 23 2 'Option D' "blue"
SetCell 24 1 $null "red"
SetCell 24 2 'Option A' "blue"
SetCell 25 1 'Option B' "red"
SetCell 25 2 'Option A' "blue"
SetCell 26 1 'Option A' "red"
SetCell 26 2 'Option C' "blue"
SetCell 27 1 'Option A' "green"
SetCell 27 2 'Option A' "blue"
SetCell 28 1 $null "red"
SetCell 28 2 'Option D' "blue"
SetCell 29 1 $null "red"
SetCell 29 2 'Option D' "blue"
SetCell 30 1 $null "red"
SetCell 30 2 'Option B' "blue"
SetCell 31 1 'Option A' "red"
SetCell 31 2 'Option D' "blue"
SetCell 32 1 'Option C' "green"
SetCell 32 2 'Option C' "blue"
SetCell 33 1 'Option D' "green"
SetCell 33 2 'Option D' "blue"
SetCell 34 1 'Option A' "red"
SetCell 34 2 'Option B' "blue"
SetCell 35 1 'Option C' "red"
SetCell 35 2 'Option D' "blue"
SetCell 36 1 'Option B' "red"
SetCell 36 2 'Option A' "blue"
SetCell 37 1 $null "red"
SetCell 37 2 'Option A' "blue"
SetCell 38 1 'Option A' "green"
SetCell 38 2 'Option A' "blue"
SetCell 39 1 $null "red"
SetCell 39 2 'Option D' "blue"
SetCell 40 1 'Option A' "red"
SetCell 40 2 'Option D' "blue"
